$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2
$ws.Range("H2").Value = 98.96429000000001
$ws.Range("I2").Value = 99.03704
$ws.Range("J2").Value = 97
$ws.Range("K2").Value = 99.03704
$ws.Range("L2").Value = 97
$ws.Range("M2").Value = 13.96296
$ws.Range("N2").Value = -323

# Row 33
$ws.Range("H33").Value = 326.8
$ws.Range("I33").Value = 350.14816
$ws.Range("J33").Value = 116.666664
$ws.Range("K33").Value = 350.14816
$ws.Range("L33").Value = 116.666664
$ws.Range("M33").Value = -121.14816
$ws.Range("N33").Value = -574.666664

# Row 129
$ws.Range("H129").Value = 20674.275
$ws.Range("I129").Value = 555.3570999999999
$ws.Range("J129").Value = 28286.838
$ws.Range("K129").Value = 1666.0713
$ws.Range("L129").Value = 84860.514
$ws.Range("M129").Value = 3333.9287
$ws.Range("N129").Value = -94860.514

# Row 132
$ws.Range("H132").Value = 1906091.5
$ws.Range("I132").Value = 1985402.1
$ws.Range("J132").Value = 2635.3333
$ws.Range("K132").Value = 5956206.300000001
$ws.Range("L132").Value = 7905.999899999999
$ws.Range("M132").Value = -5953676.300000001
$ws.Range("N132").Value = -12965.9999

# Row 141
$ws.Range("H141").Value = 1754.9231
$ws.Range("I141").Value = 1135.7567
$ws.Range("J141").Value = 3282.2
$ws.Range("K141").Value = 3407.2701
$ws.Range("L141").Value = 9846.599999999999
$ws.Range("M141").Value = 1772.7299
$ws.Range("N141").Value = -20206.6


$ws = $wb.Worksheets.Item("ARM")
# Row 4
$ws.Range("H4").Value = 115
$ws.Range("I4").Value = 115
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 115
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = 1
$ws.Range("N4").ClearContents()

# Row 5
$ws.Range("H5").Value = 41667544
$ws.Range("I5").Value = 55556224
$ws.Range("K5").Value = 55556224
$ws.Range("M5").Value = -55556112

# Row 32
$ws.Range("H32").Value = 4599.2837
$ws.Range("I32").Value = 4082.4067
$ws.Range("J32").Value = 6632.3335
$ws.Range("K32").Value = 4082.4067
$ws.Range("L32").Value = 6632.3335
$ws.Range("M32").Value = -3795.4067
$ws.Range("N32").Value = -7206.3335

# Row 45
$ws.Range("H45").Value = 1184.2632
$ws.Range("I45").Value = 992.9286
$ws.Range("J45").Value = 1720
$ws.Range("K45").Value = 992.9286
$ws.Range("L45").Value = 1720
$ws.Range("M45").Value = -615.9286
$ws.Range("N45").Value = -2474

# Row 61
$ws.Range("H61").Value = 1276.9844
$ws.Range("I61").Value = 934.9268
$ws.Range("J61").Value = 1886.7391
$ws.Range("K61").Value = 934.9268
$ws.Range("L61").Value = 1886.7391
$ws.Range("M61").Value = -722.9268
$ws.Range("N61").Value = -2310.7391

# Row 114
$ws.Range("H114").Value = 27099
$ws.Range("J114").Value = 27099
$ws.Range("L114").Value = 27099
$ws.Range("N114").Value = -35777

# Row 122
$ws.Range("H122").Value = 1147.2307
$ws.Range("I122").Value = 1158.9166
$ws.Range("J122").Value = 1007
$ws.Range("K122").Value = 3476.7498
$ws.Range("L122").Value = 3021
$ws.Range("M122").Value = -1026.7498
$ws.Range("N122").Value = -7921

# Row 136
$ws.Range("H136").Value = 1276.9844
$ws.Range("I136").Value = 934.9268
$ws.Range("J136").Value = 1886.7391
$ws.Range("K136").Value = 2804.7804
$ws.Range("L136").Value = 5660.2173
$ws.Range("M136").Value = -254.7803999999996
$ws.Range("N136").Value = -10760.2173


$ws = $wb.Worksheets.Item("BSM")
# Row 4
$ws.Range("H4").Value = 41667544
$ws.Range("I4").Value = 55556224
$ws.Range("K4").Value = 55556224
$ws.Range("M4").Value = -55556109


$ws = $wb.Worksheets.Item("CRP")
# Row 5
$ws.Range("H5").Value = 314.77777
$ws.Range("I5").Value = 219.5
$ws.Range("J5").Value = 505.33334
$ws.Range("K5").Value = 219.5
$ws.Range("L5").Value = 505.33334
$ws.Range("M5").Value = -107.5
$ws.Range("N5").Value = -729.33334

# Row 12
$ws.Range("H12").Value = 495.33334
$ws.Range("I12").Value = 495.33334
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 495.33334
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -325.33334
$ws.Range("N12").ClearContents()

# Row 19
$ws.Range("H19").Value = 620.6667
$ws.Range("I19").Value = 544.8
$ws.Range("J19").Value = 1000
$ws.Range("K19").Value = 544.8
$ws.Range("L19").Value = 1000
$ws.Range("M19").Value = -374.8
$ws.Range("N19").Value = -1340

# Row 22
$ws.Range("H22").Value = 398.2143
$ws.Range("I22").Value = 334.6842
$ws.Range("J22").Value = 532.3333
$ws.Range("K22").Value = 334.6842
$ws.Range("L22").Value = 532.3333
$ws.Range("M22").Value = 15.31580000000002
$ws.Range("N22").Value = -1232.3333

# Row 24
$ws.Range("H24").Value = 620.6667
$ws.Range("I24").Value = 544.8
$ws.Range("J24").Value = 1000
$ws.Range("K24").Value = 544.8
$ws.Range("L24").Value = 1000
$ws.Range("M24").Value = -374.8
$ws.Range("N24").Value = -1340

# Row 31
$ws.Range("H31").Value = 4507427
$ws.Range("I31").Value = 2036.75
$ws.Range("J31").Value = 12825070
$ws.Range("K31").Value = 2036.75
$ws.Range("L31").Value = 12825070
$ws.Range("M31").Value = -1741.75
$ws.Range("N31").Value = -12825660

# Row 34
$ws.Range("H34").Value = 4507427
$ws.Range("I34").Value = 2036.75
$ws.Range("J34").Value = 12825070
$ws.Range("K34").Value = 2036.75
$ws.Range("L34").Value = 12825070
$ws.Range("M34").Value = -1834.75
$ws.Range("N34").Value = -12825474


$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 175
$ws.Range("I7").Value = 100
$ws.Range("J7").Value = 250
$ws.Range("K7").Value = 300
$ws.Range("L7").Value = 750
$ws.Range("M7").Value = -188
$ws.Range("N7").Value = -974

# Row 11
$ws.Range("H11").Value = 4000725.8
$ws.Range("I11").Value = 157.3077
$ws.Range("J11").Value = 8334675
$ws.Range("K11").Value = 471.9231
$ws.Range("L11").Value = 25004025
$ws.Range("M11").Value = -331.9231
$ws.Range("N11").Value = -25004305

# Row 21
$ws.Range("H21").Value = 99
$ws.Range("I21").Value = 99
$ws.Range("K21").Value = 297
$ws.Range("M21").Value = -124

# Row 121
$ws.Range("H121").Value = 1072.6346
$ws.Range("I121").Value = 192.5
$ws.Range("J121").Value = 1232.659
$ws.Range("K121").Value = 577.5
$ws.Range("L121").Value = 3697.977
$ws.Range("M121").Value = 732.5
$ws.Range("N121").Value = -6317.977000000001

# Row 131
$ws.Range("H131").Value = 1951943.8
$ws.Range("J131").Value = 3473145.2
$ws.Range("L131").Value = 10419435.6
$ws.Range("N131").Value = -10429515.6


$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 66.666664
$ws.Range("I2").Value = 31.333334
$ws.Range("J2").Value = 102
$ws.Range("K2").Value = 31.333334
$ws.Range("L2").Value = 102
$ws.Range("M2").Value = 81.66666599999999
$ws.Range("N2").Value = -328

# Row 122
$ws.Range("H122").Value = 20410486
$ws.Range("I122").Value = 83336720
$ws.Range("J122").Value = 1978.5946
$ws.Range("K122").Value = 250010160
$ws.Range("L122").Value = 5935.783799999999
$ws.Range("M122").Value = -250007710
$ws.Range("N122").Value = -10835.7838


$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 411.90475
$ws.Range("I107").Value = 350.22223
$ws.Range("J107").Value = 458.16666
$ws.Range("K107").Value = 1050.66669
$ws.Range("L107").Value = 1374.49998
$ws.Range("M107").Value = 869.33331
$ws.Range("N107").Value = -5214.499980000001

# Row 122
$ws.Range("H122").Value = 25338.285
$ws.Range("I122").Value = 30920.5
$ws.Range("J122").Value = 1613.875
$ws.Range("K122").Value = 92761.5
$ws.Range("L122").Value = 4841.625
$ws.Range("M122").Value = -90311.5
$ws.Range("N122").Value = -9741.625

